$wb = $excel.ActiveWorkbook

# --- Rename the second sheet tab ---
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include from Prestroke and Po"

# --- Update Metadata sheet values ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.1 -> 0.0.0
$wsMeta.Range("B3").Value = "0.0.0"

# Title: drop the "ValueSet of " prefix
$wsMeta.Range("B5").Value = "Prestroke and Poststroke Functional Status: Ambulation"

# Experimental: fill in the previously empty value with the literal text "false".
# Assigning the bare string "false" via .Value gets auto-coerced to a Boolean by
# the COM layer (just like typing it into a General-formatted cell in real Excel),
# so build it indirectly: start from a padded literal (not recognized as a
# boolean keyword), trim it with a formula, then paste the computed text back in
# as a plain value. This yields a genuine shared-string "false" instead of a
# boolean cell.
$helper1 = $wsMeta.Range("ZZ1")
$helper1.Value = "false "
$helper2 = $wsMeta.Range("ZZ2")
$helper2.Formula = "=TRIM(ZZ1)"
$helper2.Copy()
$wsMeta.Range("B7").PasteSpecial(-4163)
$helper1.ClearContents()
$helper2.ClearContents()

# Date: updated timestamp
$wsMeta.Range("B8").Value = "2024-01-11T13:00:00-03:00"

# Description: reworded
$wsMeta.Range("B12").Value = "ValueSet that defines the response values for the Prestroke and Poststroke Functional Status: Ambulation."

# --- Update Include sheet values ---
$wsInclude.Range("B6").Value = "https://molic-avc.gabriellesantosleandro.com/CodeSystem/StrokeFuncStatusCS"
